# Update the cryptocurrency price/volume figures (and the InjectiveProtocol/USDe
# row re-ordering) per the scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.406.49'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '3.062.61'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'592.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').Value = "'154.31"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'0.540"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.18%  '
$ws.Range('D9').Value = '3.062.46'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').Value = "'0.155"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('D11').Value = "'5.83"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').Value = "'0.451"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.72%  '
$ws.Range('D13').Value = "'36.88"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('D14').Value = "'0.0000237"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('E15').Value = '  +1.24%  '
$ws.Range('D16').Value = '3.568.75'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').Value = "'7.17"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '63.382.34'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').Value = '3.065.43'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').Value = "'490.30"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.12%  '
$ws.Range('D21').Value = "'14.40"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.19%  '
$ws.Range('D22').Value = "'0.706"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('D23').Value = "'7.54"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').Value = "'2.44"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.65%  '
$ws.Range('D25').Value = "'81.94"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').Value = "'12.86"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('D27').Value = "'10.73"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +10.12%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = "'7.38"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.33%  '
$ws.Range('D30').Value = "'2.69"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').Value = "'2.21"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('D32').Value = "'1.00"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').Value = "'27.32"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.51%  '
$ws.Range('E34').Value = '  -1.54%  '
$ws.Range('D35').Value = "'1.06"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('D36').Value = '0.0₃0823'
$ws.Range('E36').Value = '  -3.37%  '
$ws.Range('D37').Value = "'3.33"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.09%  '
$ws.Range('D38').Value = "'5.98"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.59%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = "'9.24"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').Value = "'50.60"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').Value = "'438.71"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('D43').Value = "'0.292"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.29%  '
$ws.Range('E44').Value = '  +2.35%  '
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').Value = '2.843.69'
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('D47').Value = "'38.94"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('D48').Value = "'130.17"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').Value = "'1.00"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = "'25.23"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('E51').Value = '  -1.42%  '
